$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 639.9091
$ws.Range("J17").Value = 639.9091
$ws.Range("L17").Value = 1919.7273
$ws.Range("N17").Value = -2255.7273
$ws.Range("H38").Value = 3349.5
$ws.Range("J38").Value = 6499
$ws.Range("L38").Value = 19497
$ws.Range("N38").Value = -20241
$ws.Range("H51").Value = 1300
$ws.Range("J51").Value = 1300
$ws.Range("L51").Value = 1300
$ws.Range("N51").Value = -2268
$ws.Range("H87").Value = 39986
$ws.Range("J87").Value = 39986
$ws.Range("L87").Value = 39986
$ws.Range("N87").Value = -42482
$ws.Range("H90").Value = 39986
$ws.Range("J90").Value = 39986
$ws.Range("L90").Value = 119958
$ws.Range("N90").Value = -132438
$ws.Range("H92").Value = 2466.6667
$ws.Range("I92").Value = 2466.6667
$ws.Range("K92").Value = 2466.6667
$ws.Range("M92").Value = -1218.6667
$ws.Range("H116").Value = 5064.5386
$ws.Range("I116").Value = 4993.3335
$ws.Range("K116").Value = 4993.3335
$ws.Range("M116").Value = -1551.3335
$ws.Range("H125").Value = 6629.3335
$ws.Range("I125").Value = 6899
$ws.Range("J125").Value = 6494.5
$ws.Range("K125").Value = 62091
$ws.Range("L125").Value = 58450.5
$ws.Range("M125").Value = -59631
$ws.Range("N125").Value = -63370.5
$ws.Range("H135").Value = 6023.7144
$ws.Range("J135").Value = 1418
$ws.Range("L135").Value = 12762
$ws.Range("N135").Value = -17832
$ws.Range("H138").Value = 8030.647
$ws.Range("J138").Value = 8157.5625
$ws.Range("L138").Value = 24472.6875
$ws.Range("N138").Value = -34752.6875
$ws.Range("H141").Value = 4500.5
$ws.Range("I141").Value = 3417.3333
$ws.Range("K141").Value = 10251.9999
$ws.Range("M141").Value = -5071.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 5000
$ws.Range("K30").Value = 5000
$ws.Range("M30").Value = -4850
$ws.Range("H32").Value = 4105.55
$ws.Range("I32").Value = 4111.1055
$ws.Range("K32").Value = 4111.1055
$ws.Range("M32").Value = -3824.1055
$ws.Range("H88").Value = 999
$ws.Range("I88").Value = 999
$ws.Range("K88").Value = 999
$ws.Range("M88").Value = -593
$ws.Range("H91").Value = 999
$ws.Range("I91").Value = 999
$ws.Range("K91").Value = 999
$ws.Range("M91").Value = 405

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 734.1875
$ws.Range("I94").Value = 587.25
$ws.Range("J94").Value = 1175
$ws.Range("K94").Value = 587.25
$ws.Range("L94").Value = 1175
$ws.Range("M94").Value = -136.25
$ws.Range("N94").Value = -2077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9890.1
$ws.Range("I31").Value = 2200
$ws.Range("K31").Value = 2200
$ws.Range("M31").Value = -1905
$ws.Range("H34").Value = 9890.1
$ws.Range("I34").Value = 2200
$ws.Range("K34").Value = 2200
$ws.Range("M34").Value = -1998
$ws.Range("H58").Value = 2278.2
$ws.Range("I58").Value = 2123
$ws.Range("J58").Value = 2899
$ws.Range("K58").Value = 2123
$ws.Range("L58").Value = 2899
$ws.Range("M58").Value = -1920
$ws.Range("N58").Value = -3305
$ws.Range("H74").Value = 36307.43
$ws.Range("J74").Value = 36307.43
$ws.Range("L74").Value = 36307.43
$ws.Range("N74").Value = -38055.43
$ws.Range("H77").Value = 36307.43
$ws.Range("J77").Value = 36307.43
$ws.Range("L77").Value = 108922.29
$ws.Range("N77").Value = -117658.29
$ws.Range("H99").Value = 1999
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1999
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1999
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4995
$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1999
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5997
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -10937
$ws.Range("H132").Value = 2581.7
$ws.Range("I132").Value = 1304.0834
$ws.Range("K132").Value = 3912.2502
$ws.Range("M132").Value = -1382.2502
$ws.Range("H134").Value = 2486.3076
$ws.Range("I134").Value = 2360.1667
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 7080.500100000001
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4545.500100000001
$ws.Range("N134").Value = -17070
$ws.Range("H136").Value = 2278.2
$ws.Range("I136").Value = 2123
$ws.Range("J136").Value = 2899
$ws.Range("K136").Value = 6369
$ws.Range("L136").Value = 8697
$ws.Range("M136").Value = -3819
$ws.Range("N136").Value = -13797

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 217.16667
$ws.Range("I38").Value = 187.75
$ws.Range("K38").Value = 563.25
$ws.Range("M38").Value = -216.25
$ws.Range("H132").Value = 3952.8
$ws.Range("J132").Value = 5612
$ws.Range("L132").Value = 50508
$ws.Range("N132").Value = -55568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 11697
$ws.Range("J96").Value = 11697
$ws.Range("L96").Value = 11697
$ws.Range("N96").Value = -17189
$ws.Range("H132").Value = 3002.2
$ws.Range("I132").Value = 2484.111
$ws.Range("K132").Value = 7452.333
$ws.Range("M132").Value = -4922.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H40").Value = 3235.3125
$ws.Range("I40").Value = 3183.7856
$ws.Range("J40").Value = 3596
$ws.Range("K40").Value = 3183.7856
$ws.Range("L40").Value = 3596
$ws.Range("M40").Value = -3047.7856
$ws.Range("N40").Value = -3868
$ws.Range("H55").Value = 666.1
$ws.Range("J55").Value = 945
$ws.Range("L55").Value = 945
$ws.Range("N55").Value = -1291
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3114.4211
$ws.Range("I132").Value = 2345.4
$ws.Range("K132").Value = 7036.200000000001
$ws.Range("M132").Value = -4506.200000000001
$ws.Range("H136").Value = 35187.445
$ws.Range("J136").Value = 48331.25
$ws.Range("L136").Value = 144993.75
$ws.Range("N136").Value = -150093.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3317
$ws.Range("I96").Value = 3473.25
$ws.Range("J96").Value = 2067
$ws.Range("K96").Value = 3473.25
$ws.Range("L96").Value = 2067
$ws.Range("M96").Value = -2100.25
$ws.Range("N96").Value = -4813
$ws.Range("H100").Value = 795.6
$ws.Range("J100").Value = 869.75
$ws.Range("L100").Value = 1739.5
$ws.Range("N100").Value = -2821.5
$ws.Range("H122").Value = 1394.4445
$ws.Range("I122").Value = 1241.2354
$ws.Range("K122").Value = 3723.7062
$ws.Range("M122").Value = -1273.7062
$ws.Range("H126").Value = 1967.8334
$ws.Range("I126").Value = 1967.8334
$ws.Range("K126").Value = 5903.5002
$ws.Range("M126").Value = -3433.5002
$ws.Range("H132").Value = 2036.7073
$ws.Range("I132").Value = 1810.3549
$ws.Range("K132").Value = 5431.0647
$ws.Range("M132").Value = -2901.0647
$ws.Range("H136").Value = 7128.095
$ws.Range("I136").Value = 9238.333000000001
$ws.Range("J136").Value = 1852.5
$ws.Range("K136").Value = 27714.999
$ws.Range("L136").Value = 5557.5
$ws.Range("M136").Value = -25164.999
$ws.Range("N136").Value = -10657.5

Write-Output "applied changes"
